$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: change phone number "0909099999" (text) to numeric 909099999
$ws.Range("C2").Value = 909099999

# G2: set membership card number "111111111" as text (leading apostrophe forces text entry, like typing into Excel)
$ws.Range("G2").Value = "'111111111"
$ws.Range("G2").Style = "Normal"
